$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.556.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.05%  "
$ws.Range("D3").Value = "'1.599.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "'212.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'0.513"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'26.82"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.54%  "
$ws.Range("D9").Value = "'43.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.17%  "
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'1.825.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "'1.597.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "'29.560.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("E16").Value = "  +3.48%  "
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "'63.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("D19").Value = "'239.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D23").Value = "'3.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").Value = "'9.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "'154.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("E27").Value = "  +3.25%  "
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'3.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.21%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "'1.435.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'1.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").Value = "'0.539"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.0494"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.09%  "
$ws.Range("D44").Value = "'53.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +26.34%  "
$ws.Range("D45").Value = "'0.801"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.23%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "'0.989"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +19.40%  "
$ws.Range("D48").Value = "'65.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.29%  "
$ws.Range("D49").Value = "'5.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").Value = "'1.737.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").Value = "'86.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.71%  "
